# Auto-generated PowerShell Excel COM-interop script
# Applies the Cerberus_Profits.xlsx market-data refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 507.33334
$ws.Range("I12").Value = 154.33333
$ws.Range("J12").Value = 1566.3334
$ws.Range("K12").Value = 154.33333
$ws.Range("L12").Value = 1566.3334
$ws.Range("M12").Value = 15.66667000000001
$ws.Range("N12").Value = -1906.3334
$ws.Range("H51").Value = 9762.182000000001
$ws.Range("I51").Value = 14249.25
$ws.Range("K51").Value = 14249.25
$ws.Range("M51").Value = -13765.25
$ws.Range("H74").Value = 5849.2593
$ws.Range("I74").Value = 4859.5713
$ws.Range("K74").Value = 4859.5713
$ws.Range("M74").Value = -3923.5713
$ws.Range("H77").Value = 5849.2593
$ws.Range("I77").Value = 4859.5713
$ws.Range("K77").Value = 24297.8565
$ws.Range("M77").Value = -19617.8565
$ws.Range("H101").Value = 2882.353
$ws.Range("I101").Value = 1391.8334
$ws.Range("K101").Value = 4175.5002
$ws.Range("M101").Value = -2553.5002
$ws.Range("H113").Value = 7035.1113
$ws.Range("I113").Value = 6067.357
$ws.Range("K113").Value = 6067.357
$ws.Range("M113").Value = -2813.357
$ws.Range("H116").Value = 6540
$ws.Range("I116").Value = 9814.666999999999
$ws.Range("K116").Value = 9814.666999999999
$ws.Range("M116").Value = -6372.666999999999
$ws.Range("H132").Value = 3786.6667
$ws.Range("I132").Value = 3701.7659
$ws.Range("K132").Value = 11105.2977
$ws.Range("M132").Value = -8575.297699999999
$ws.Range("H138").Value = 2737.7673
$ws.Range("I138").Value = 2228.7368
$ws.Range("K138").Value = 6686.2104
$ws.Range("M138").Value = -1546.2104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 739.2727
$ws.Range("I2").Value = 336.16666
$ws.Range("J2").Value = 1223
$ws.Range("K2").Value = 336.16666
$ws.Range("L2").Value = 1223
$ws.Range("M2").Value = -223.16666
$ws.Range("N2").Value = -1449
$ws.Range("H16").Value = 21332.666
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 21332.666
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 21332.666
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -21906.666
$ws.Range("H32").Value = 3139.131
$ws.Range("I32").Value = 2278.7593
$ws.Range("K32").Value = 2278.7593
$ws.Range("M32").Value = -1991.7593
$ws.Range("H45").Value = 2070
$ws.Range("I45").Value = 1583.3334
$ws.Range("J45").Value = 2232.2222
$ws.Range("K45").Value = 1583.3334
$ws.Range("L45").Value = 2232.2222
$ws.Range("M45").Value = -1206.3334
$ws.Range("N45").Value = -2986.2222
$ws.Range("H116").Value = 739.2727
$ws.Range("I116").Value = 336.16666
$ws.Range("J116").Value = 1223
$ws.Range("K116").Value = 336.16666
$ws.Range("L116").Value = 1223
$ws.Range("M116").Value = 1957.83334
$ws.Range("N116").Value = -5811
$ws.Range("H132").Value = 1984.1136
$ws.Range("I132").Value = 1823.5581
$ws.Range("K132").Value = 5470.6743
$ws.Range("M132").Value = -2940.6743

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 739.2727
$ws.Range("I3").Value = 336.16666
$ws.Range("J3").Value = 1223
$ws.Range("K3").Value = 336.16666
$ws.Range("L3").Value = 1223
$ws.Range("M3").Value = -222.16666
$ws.Range("N3").Value = -1451
$ws.Range("H132").Value = 91000
$ws.Range("J132").Value = 91000
$ws.Range("L132").Value = 91000
$ws.Range("N132").Value = -101120
$ws.Range("H133").Value = 79244.75
$ws.Range("J133").Value = 85659.664
$ws.Range("L133").Value = 85659.664
$ws.Range("N133").Value = -95779.664
$ws.Range("H134").Value = 5576.589
$ws.Range("I134").Value = 5232.456
$ws.Range("K134").Value = 15697.368
$ws.Range("M134").Value = -13162.368
$ws.Range("H139").Value = 197799.6
$ws.Range("J139").Value = 197799.6
$ws.Range("L139").Value = 197799.6
$ws.Range("N139").Value = -208079.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4067
$ws.Range("I31").Value = 2144.25
$ws.Range("K31").Value = 2144.25
$ws.Range("M31").Value = -1849.25
$ws.Range("H34").Value = 4067
$ws.Range("I34").Value = 2144.25
$ws.Range("K34").Value = 2144.25
$ws.Range("M34").Value = -1942.25
$ws.Range("H58").Value = 1076.3024
$ws.Range("I58").Value = 676.6286
$ws.Range("K58").Value = 676.6286
$ws.Range("M58").Value = -473.6286
$ws.Range("H98").Value = 16780
$ws.Range("J98").Value = 16780
$ws.Range("L98").Value = 16780
$ws.Range("N98").Value = -21272
$ws.Range("H132").Value = 2276.2622
$ws.Range("I132").Value = 2331.8928
$ws.Range("K132").Value = 6995.678400000001
$ws.Range("M132").Value = -4465.678400000001
$ws.Range("H136").Value = 1076.3024
$ws.Range("I136").Value = 676.6286
$ws.Range("K136").Value = 2029.8858
$ws.Range("M136").Value = 520.1142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 139038660
$ws.Range("I4").Value = 185134880
$ws.Range("J4").Value = 750000
$ws.Range("K4").Value = 555404640
$ws.Range("L4").Value = 2250000
$ws.Range("M4").Value = -555404528
$ws.Range("N4").Value = -2250224
$ws.Range("H23").Value = 648.35297
$ws.Range("I23").Value = 469.4
$ws.Range("J23").Value = 722.9167
$ws.Range("K23").Value = 1408.2
$ws.Range("L23").Value = 2168.7501
$ws.Range("M23").Value = -1173.2
$ws.Range("N23").Value = -2638.7501
$ws.Range("H56").Value = 6791.6343
$ws.Range("I56").Value = 6791.6343
$ws.Range("K56").Value = 6791.6343
$ws.Range("M56").Value = -6261.6343
$ws.Range("H104").Value = 4945.1665
$ws.Range("I104").Value = 3495
$ws.Range("J104").Value = 5235.2
$ws.Range("K104").Value = 10485
$ws.Range("L104").Value = 15705.6
$ws.Range("M104").Value = -7864
$ws.Range("N104").Value = -20947.6
$ws.Range("H132").Value = 3559.3333
$ws.Range("I132").Value = 3439.5
$ws.Range("J132").Value = 3799
$ws.Range("K132").Value = 30955.5
$ws.Range("L132").Value = 34191
$ws.Range("M132").Value = -28425.5
$ws.Range("N132").Value = -39251
$ws.Range("H139").Value = 13895950
$ws.Range("I139").Value = 16669140
$ws.Range("J139").Value = 30000
$ws.Range("K139").Value = 50007420
$ws.Range("L139").Value = 90000
$ws.Range("M139").Value = -50002280
$ws.Range("N139").Value = -100280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H102").Value = 3497.1143
$ws.Range("I102").Value = 3450.2222
$ws.Range("K102").Value = 3450.2222
$ws.Range("M102").Value = -1828.2222
$ws.Range("H113").Value = 2303.5833
$ws.Range("I113").Value = 2210.5
$ws.Range("K113").Value = 2210.5
$ws.Range("M113").Value = -40.5
$ws.Range("H122").Value = 2890.8262
$ws.Range("I122").Value = 2408.9524
$ws.Range("J122").Value = 7950.5
$ws.Range("K122").Value = 7226.8572
$ws.Range("L122").Value = 23851.5
$ws.Range("M122").Value = -4776.8572
$ws.Range("N122").Value = -28751.5
$ws.Range("H132").Value = 2142.1553
$ws.Range("I132").Value = 2086.4565
$ws.Range("J132").Value = 2355.6667
$ws.Range("K132").Value = 6259.369499999999
$ws.Range("L132").Value = 7067.000100000001
$ws.Range("M132").Value = -3729.369499999999
$ws.Range("N132").Value = -12127.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2289.2
$ws.Range("I46").Value = 1360
$ws.Range("J46").Value = 3218.4
$ws.Range("K46").Value = 1360
$ws.Range("L46").Value = 3218.4
$ws.Range("M46").Value = -1172
$ws.Range("N46").Value = -3594.4
$ws.Range("H61").Value = 2412.875
$ws.Range("I61").Value = 2717.1667
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 2717.1667
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -2515.1667
$ws.Range("N61").Value = -1904
$ws.Range("H93").Value = 1118.0625
$ws.Range("I93").Value = 1008
$ws.Range("J93").Value = 1301.5
$ws.Range("K93").Value = 1008
$ws.Range("L93").Value = 1301.5
$ws.Range("M93").Value = 240
$ws.Range("N93").Value = -3797.5
$ws.Range("H113").Value = 2412.875
$ws.Range("I113").Value = 2717.1667
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 2717.1667
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -547.1667000000002
$ws.Range("N113").Value = -5840
$ws.Range("H132").Value = 1939.8939
$ws.Range("I132").Value = 1705.84
$ws.Range("K132").Value = 5117.52
$ws.Range("M132").Value = -2587.52

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4280.085
$ws.Range("I132").Value = 4558.615
$ws.Range("J132").Value = 2211
$ws.Range("K132").Value = 13675.845
$ws.Range("L132").Value = 6633
$ws.Range("M132").Value = -11145.845
$ws.Range("N132").Value = -11693
